# 自动更新Excel文件 - 2025-10-23 23:11:37
# This script advances the "剩余" (remaining days) counter for every data row
# by one day, as if a day has passed since the sheet was last generated.
# Column layout: A=行号 B=店铺名称 C=地址 D=总天(total days) E=剩余(remaining days)
#                F=开始时间(start date, yyyymmdd) G=备注1 H=备注2 I=备注3
#
# Rule observed from the source data:
#   - For every data row, E (remaining) is decreased by 1.
#   - If that would bring E to 0 (i.e. the cycle expired), instead the row is
#     "renewed": E is reset back to the row's total (D), and the start date F
#     is pushed forward by D days (so the cycle restarts as of the new day).
#   - Rows whose start date (F) is not a valid yyyymmdd date are left
#     untouched (data entry anomaly, e.g. row 36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($null -eq $dVal -or $null -eq $eVal -or $null -eq $fVal) {
        continue
    }

    # Skip rows whose start date isn't a well-formed 8-digit yyyymmdd value
    $fStr = [string]([int]$fVal)
    if ($fStr.Length -ne 8) {
        continue
    }

    $newE = [int]$eVal - 1
    if ($newE -le 0) {
        # Cycle expired: renew it starting today, full duration restored.
        $eCell.Value2 = [int]$dVal
        $fCell.Value2 = [int]$fVal + [int]$dVal
    } else {
        $eCell.Value2 = $newE
    }
}
